$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to make room for today's price entry.
# This shifts all existing rows (2..198) down by one (to 3..199),
# preserving their values, number formats and hyperlinks' cell anchors.
$ws.Rows.Item(2).Insert()

# Copy the formatting from the row below (the row that used to be row 2,
# now shifted to row 3) onto the newly inserted row 2, so the new row
# starts with the same look (centered text / number style) as the rest
# of the table instead of inheriting the header row's format.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Force the text columns to be stored as plain text so that date-like
# strings (e.g. "20-02-2026") are not auto-converted into date serials.
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("E2:F2").NumberFormat = "@"

# Populate the new top row with the newest price entry.
$ws.Range("A2").Value = "20-02-2026"
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 320.45
$ws.Range("E2").Value = "01-02-2026"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-02-2026.pdf"

# Give the new row's Circular Link cell its own hyperlink.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-02-2026.pdf")

# The bottom-most historical row (formerly row 198) was pushed down to
# row 199 by the insert above; give its link cell a hyperlink too (it
# duplicates the previous last row's data/link).
$ws.Hyperlinks.Add($ws.Range("F199"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# Inserting a row shifts cell values/content down correctly, but the
# existing hyperlink objects stay anchored to their original cell
# addresses (their underlying link target does not travel with the
# shifted text). Re-point every link cell so that its hyperlink target
# matches the (now shifted) displayed URL text in that cell.
for ($r = 3; $r -le 198; $r++) {
    $cell = $ws.Range("F$r")
    $url = $cell.Value
    $existing = $cell.Hyperlinks.Item(1)
    if ($existing.Address -ne $url) {
        $ws.Hyperlinks.Add($cell, $url)
    }
}
